$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its exact text formatting (e.g. trailing
# zeros, thousand-separator dots) instead of being auto-converted to a number.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '71.321.50'
$ws.Range('E2').Value = '  +2.69%  '
$ws.Range('D3').Value = '3.701.04'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '582.78'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = '177.93'
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('D7').Value = '3.691.39'
$ws.Range('E7').Value = '  +8.15%  '
$ws.Range('E8').Value = '  +4.11%  '
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').Value = '0.200'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').Value = '6.82'
$ws.Range('E11').Value = '  +26.50%  '
$ws.Range('D12').Value = '0.610'
$ws.Range('E12').Value = '  +4.91%  '
$ws.Range('D13').Value = '49.11'
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '4.293.36'
$ws.Range('E15').Value = '  +8.23%  '
$ws.Range('D16').Value = '678.05'
$ws.Range('E16').Value = '  -2.57%  '
$ws.Range('D17').Value = '9.00'
$ws.Range('E17').Value = '  +4.67%  '
$ws.Range('D18').Value = '3.697.95'
$ws.Range('E18').Value = '  +8.19%  '
$ws.Range('D19').Value = '71.435.33'
$ws.Range('E19').Value = '  +2.83%  '
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').Value = '17.98'
$ws.Range('E21').Value = '  +2.09%  '
$ws.Range('D22').Value = '11.60'
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('D23').Value = '0.942'
$ws.Range('E23').Value = '  +5.49%  '
$ws.Range('D24').Value = '17.44'
$ws.Range('E24').Value = '  +3.29%  '
$ws.Range('D25').Value = '102.09'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').Value = '3.98'
$ws.Range('E26').Value = '  +2.18%  '
$ws.Range('D27').Value = '2.83'
$ws.Range('E27').Value = '  +6.63%  '
$ws.Range('D28').Value = '10.31'
$ws.Range('E28').Value = '  +8.25%  '
$ws.Range('D29').Value = '35.21'
$ws.Range('E29').Value = '  +5.30%  '
$ws.Range('B30').Value = 'Stacks'
$ws.Range('C30').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D30').Value = '3.43'
$ws.Range('E30').Value = '  +5.46%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '9.20'
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('B32').Value = 'Mantle'
$ws.Range('C32').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D32').Value = '1.43'
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = '7.52'
$ws.Range('E33').Value = '  +6.96%  '
$ws.Range('B34').Value = 'dogwifhat'
$ws.Range('C34').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D34').Value = '4.08'
$ws.Range('E34').Value = '  +10.49%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').Value = '585.77'
$ws.Range('E35').Value = '  +1.57%  '
$ws.Range('B36').Value = 'Cosmos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D36').Value = '11.18'
$ws.Range('E36').Value = '  +1.77%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.109'
$ws.Range('E37').Value = '  +5.82%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = '58.86'
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '3.681.85'
$ws.Range('E40').Value = '  +4.16%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.145'
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('D43').Value = '0.352'
$ws.Range('E43').Value = '  +6.33%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0456'
$ws.Range('E44').Value = '  +9.65%  '
$ws.Range('D45').Value = '0.0₃0766'
$ws.Range('E45').Value = '  +5.15%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '35.65'
$ws.Range('E46').Value = '  +2.54%  '
$ws.Range('D47').Value = '2.77'
$ws.Range('E47').Value = '  +4.39%  '
$ws.Range('D48').Value = '2.90'
$ws.Range('E48').Value = '  +10.08%  '
$ws.Range('E49').Value = '  +4.02%  '
$ws.Range('D50').Value = '135.62'
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('E51').Value = '  +10.41%  '
